$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# followed by the 16 numeric TPM-derived columns E..T) replacing the old rows 2-7
# with 3 recomputed rows.

$data = @(
  @("ECs",  "Guca2a", "Gucy2c", "Resolving-Mac", 1, 0.3333333333333333, 0.486124, 1.458372, 0.2724081688156144, 0.2724081688156144, 3, 1, 0.2121096666666666, 0.6363289999999999, 1, 1, 0.1031115995986667, 0.9280043963879999, 0.2724081688156144, 0.2724081688156144),
  @("FAPs", "Guca2a", "Gucy2c", "Resolving-Mac", 3, 1, 0.765631, 2.296893, 0.4290348526270409, 0.4290348526270409, 3, 1, 0.2121096666666666, 0.6363289999999999, 1, 1, 0.1623977361996666, 1.461579625797, 0.4290348526270409, 0.4290348526270409),
  @("MuSCs","Guca2a", "Gucy2c", "Resolving-Mac", 3, 1, 0.5327876666666667, 1.598363, 0.2985569785573446, 0.2985569785573447, 3, 1, 0.2121096666666666, 0.6363289999999999, 1, 1, 0.1130094143807778, 1.017084729427, 0.2985569785573446, 0.2985569785573447)
)

# Remove the three extra old rows (rows 5,6,7 -> only rows 2,3,4 remain after the edit).
$ws.Range("A5:A7").EntireRow.Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
